$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the new "Companero" / "Probabilidad" columns (B, C).
#    Probability values are stored as literal TEXT ("0.1" and not the number
#    0.1), so we enter them with a leading apostrophe to force text, matching
#    the shared-string layout of the target workbook.
# ---------------------------------------------------------------------------

$ws.Range("B1").Value = "Compañero"
$ws.Range("C1").Value = "Probabilidad"

$ws.Range("B2").Value = "Ivan"
$ws.Range("C2").Value = "'0.1"

$ws.Range("B3").Value = "Juan"
$ws.Range("C3").Value = "'0.1"

$ws.Range("B4").Value = "William"
$ws.Range("C4").Value = "'0.7"

$ws.Range("B5").Value = "Ruben"
$ws.Range("C5").Value = "'0.5"

$ws.Range("B6").Value = "Juan"
$ws.Range("C6").Value = "'0.4"

$ws.Range("B7").Value = "William"
$ws.Range("C7").Value = "'0.3"

# ---------------------------------------------------------------------------
# 2) Column widths for the new columns.
# ---------------------------------------------------------------------------

$ws.Columns("B").ColumnWidth = 19.140625
$ws.Columns("C").ColumnWidth = 20.28515625
$ws.Columns("D").ColumnWidth = 9.140625

# ---------------------------------------------------------------------------
# 3) The new cells (B1:C7, plus the formatted-but-empty D1:D9/H9 helper
#    cells below) reuse the workbook's pre-existing alternate cell style
#    (style index 1) instead of the default one, exactly like the lone
#    G8 cell already did in the source workbook. We copy that format
#    across with Copy/PasteSpecial (format-only) so the *same* style slot
#    is reused rather than a brand new one being allocated, then drop the
#    underline that style used to carry (the author also stripped the
#    now-unused underlined font from styles.xml).
# ---------------------------------------------------------------------------

$ws.Range("G8").Copy()
$ws.Range("B1:D9").PasteSpecial(-4122)
$ws.Range("H9").Select()
$ws.Range("G8").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B1:D9").Font.Underline = $false
$ws.Range("G8").Font.Underline = $false
$ws.Range("H9").Font.Underline = $false

# Remove the leftover placeholder content in G8/H9 (kept only to extend the
# sheet's dimension/used-range down to row 9 and across to column H, as in
# the target file) while leaving every other style slot untouched.
$ws.Range("G8").ClearContents()
$ws.Range("H9").ClearContents()

# ---------------------------------------------------------------------------
# 4) Final selection / view state, matching the committed workbook.
# ---------------------------------------------------------------------------

$ws.Range("F9").Select()
